$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string value while forcing Text format
# so Excel does not silently convert it to a floating point number and
# drop significant trailing zeros (e.g. "5.260" -> 5.26).
function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue $ws "D2" "243.79"
Set-TextValue $ws "D4" "5.260"
Set-TextValue $ws "D5" "0.05842"
Set-TextValue $ws "D6" "6.457"
Set-TextValue $ws "D7" "3.337"
Set-TextValue $ws "D8" "0.8078"
Set-TextValue $ws "D9" "0.8999"
Set-TextValue $ws "D10" "0.1380"
Set-TextValue $ws "D11" "0.07080"
Set-TextValue $ws "D12" "0.03067"
Set-TextValue $ws "D13" "0.03028"
Set-TextValue $ws "D14" "0.09331"
Set-TextValue $ws "D15" "3.813"
Set-TextValue $ws "D16" "0.001544"
Set-TextValue $ws "D17" "0.04718"
Set-TextValue $ws "D18" "0.0006020"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue $ws "D19" "0.006159"
Set-TextValue $ws "D20" "0.001259"
Set-TextValue $ws "D21" "0.004057"
Set-TextValue $ws "D22" "0.00008697"
Set-TextValue $ws "D24" "2.180"
Set-TextValue $ws "D25" "0.3170"
Set-TextValue $ws "D26" "0.1318"
Set-TextValue $ws "D28" "0.0002328"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006248"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1052"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002526"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws "D44" "0.006939"
Set-TextValue $ws "D45" "0.00005338"
Set-TextValue $ws "D47" "0.5100"
Set-TextValue $ws "D48" "0.007740"
